$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 14900
$ws.Range("J3").Value = 14900
$ws.Range("L3").Value = 14900
$ws.Range("N3").Value = -15128
$ws.Range("H9").Value = 125.75
$ws.Range("I9").Value = 107.5
$ws.Range("J9").Value = 144
$ws.Range("K9").Value = 107.5
$ws.Range("L9").Value = 144
$ws.Range("M9").Value = 61.5
$ws.Range("N9").Value = -482
$ws.Range("H17").Value = 1183764.8
$ws.Range("I17").Value = 2936.25
$ws.Range("J17").Value = 1243553.5
$ws.Range("K17").Value = 8808.75
$ws.Range("L17").Value = 3730660.5
$ws.Range("M17").Value = -8640.75
$ws.Range("N17").Value = -3730996.5
$ws.Range("H19").Value = 1040.5714
$ws.Range("I19").Value = 1031.1111
$ws.Range("J19").Value = 1057.6
$ws.Range("K19").Value = 1031.1111
$ws.Range("L19").Value = 1057.6
$ws.Range("M19").Value = -856.1111000000001
$ws.Range("N19").Value = -1407.6
$ws.Range("H32").Value = 48571.285
$ws.Range("I32").Value = 73500
$ws.Range("J32").Value = 15333
$ws.Range("K32").Value = 73500
$ws.Range("L32").Value = 15333
$ws.Range("M32").Value = -73174
$ws.Range("N32").Value = -15985
$ws.Range("H34").Value = 6850.1113
$ws.Range("I34").Value = 6850.1113
$ws.Range("K34").Value = 6850.1113
$ws.Range("M34").Value = -6647.1113
$ws.Range("H36").Value = 6850.1113
$ws.Range("I36").Value = 6850.1113
$ws.Range("K36").Value = 6850.1113
$ws.Range("M36").Value = -6135.1113
$ws.Range("H37").Value = 5000
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15252
$ws.Range("H43").Value = 4924.6665
$ws.Range("H62").Value = 4242.3335
$ws.Range("I62").Value = 4030.9
$ws.Range("J62").Value = 4665.2
$ws.Range("K62").Value = 4030.9
$ws.Range("L62").Value = 4665.2
$ws.Range("M62").Value = -3406.9
$ws.Range("N62").Value = -5913.2
$ws.Range("H65").Value = 4242.3335
$ws.Range("I65").Value = 4030.9
$ws.Range("J65").Value = 4665.2
$ws.Range("K65").Value = 20154.5
$ws.Range("L65").Value = 23326
$ws.Range("M65").Value = -17034.5
$ws.Range("N65").Value = -29566
$ws.Range("H74").Value = 3548.625
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 3548.625
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 2894.05
$ws.Range("I86").Value = 1917.3125
$ws.Range("K86").Value = 1917.3125
$ws.Range("M86").Value = -794.3125
$ws.Range("H89").Value = 2894.05
$ws.Range("I89").Value = 1917.3125
$ws.Range("K89").Value = 9586.5625
$ws.Range("M89").Value = -3970.5625
$ws.Range("H102").Value = 14900
$ws.Range("J102").Value = 14900
$ws.Range("L102").Value = 14900
$ws.Range("N102").Value = -21390
$ws.Range("H106").Value = 22943.15
$ws.Range("I106").Value = 6431.1113
$ws.Range("K106").Value = 6431.1113
$ws.Range("M106").Value = -5800.1113
$ws.Range("H132").Value = 2802
$ws.Range("I132").Value = 2331.4211
$ws.Range("K132").Value = 6994.263300000001
$ws.Range("M132").Value = -4464.263300000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = None
$ws.Range("N133").ClearContents()
$ws.Range("H138").Value = 3331.879
$ws.Range("I138").Value = 3125.7144
$ws.Range("J138").Value = 3692.6667
$ws.Range("K138").Value = 9377.143199999999
$ws.Range("L138").Value = 11078.0001
$ws.Range("M138").Value = -4237.143199999999
$ws.Range("N138").Value = -21358.0001
$ws.Range("H141").Value = 2584.9678
$ws.Range("I141").Value = 2540.5386
$ws.Range("K141").Value = 7621.6158
$ws.Range("M141").Value = -2441.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 9999
$ws.Range("I29").Value = 9999
$ws.Range("K29").Value = 9999
$ws.Range("M29").Value = -9691
$ws.Range("H32").Value = 5004.5293
$ws.Range("I32").Value = 3044.7307
$ws.Range("K32").Value = 3044.7307
$ws.Range("M32").Value = -2757.7307
$ws.Range("H45").Value = 7066.727
$ws.Range("I45").Value = 7563
$ws.Range("K45").Value = 7563
$ws.Range("M45").Value = -7186
$ws.Range("H61").Value = 5312.6
$ws.Range("I61").Value = 3252.7407
$ws.Range("K61").Value = 3252.7407
$ws.Range("M61").Value = -3040.7407
$ws.Range("H74").Value = 2643.1162
$ws.Range("I74").Value = 2306.2307
$ws.Range("J74").Value = 3158.353
$ws.Range("K74").Value = 2306.2307
$ws.Range("L74").Value = 3158.353
$ws.Range("M74").Value = -1432.2307
$ws.Range("N74").Value = -4906.353
$ws.Range("H77").Value = 2643.1162
$ws.Range("I77").Value = 2306.2307
$ws.Range("J77").Value = 3158.353
$ws.Range("K77").Value = 11531.1535
$ws.Range("L77").Value = 15791.765
$ws.Range("M77").Value = -7163.1535
$ws.Range("N77").Value = -24527.765
$ws.Range("H102").Value = 3888.4285
$ws.Range("I102").Value = 3969.8333
$ws.Range("K102").Value = 3969.8333
$ws.Range("M102").Value = -2347.8333
$ws.Range("H110").Value = 1558.5
$ws.Range("I110").Value = 1182.3
$ws.Range("K110").Value = 1182.3
$ws.Range("M110").Value = 862.7
$ws.Range("H112").Value = 10129
$ws.Range("J112").Value = 10129
$ws.Range("L112").Value = 10129
$ws.Range("N112").Value = -13083
$ws.Range("H122").Value = 1569.5625
$ws.Range("I122").Value = 1401.1666
$ws.Range("K122").Value = 4203.4998
$ws.Range("M122").Value = -1753.4998
$ws.Range("H124").Value = 29799.6
$ws.Range("J124").Value = 29799.6
$ws.Range("L124").Value = 29799.6
$ws.Range("N124").Value = -39619.6
$ws.Range("H132").Value = 2497.4827
$ws.Range("I132").Value = 2569.037
$ws.Range("K132").Value = 7707.110999999999
$ws.Range("M132").Value = -5177.110999999999
$ws.Range("H136").Value = 5312.6
$ws.Range("I136").Value = 3252.7407
$ws.Range("K136").Value = 9758.222099999999
$ws.Range("M136").Value = -7208.222099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 15590
$ws.Range("J88").Value = 15590
$ws.Range("L88").Value = 15590
$ws.Range("N88").Value = -16402
$ws.Range("H91").Value = 15590
$ws.Range("J91").Value = 15590
$ws.Range("L91").Value = 15590
$ws.Range("N91").Value = -18398
$ws.Range("H94").Value = 1124.75
$ws.Range("I94").Value = 854.5
$ws.Range("J94").Value = 1395
$ws.Range("K94").Value = 854.5
$ws.Range("L94").Value = 1395
$ws.Range("M94").Value = -403.5
$ws.Range("N94").Value = -2297
$ws.Range("H96").Value = 7935.3335
$ws.Range("I96").Value = 7935.3335
$ws.Range("K96").Value = 7935.3335
$ws.Range("M96").Value = -5189.3335
$ws.Range("H107").Value = 1348.5238
$ws.Range("I107").Value = 1405.4667
$ws.Range("J107").Value = 1206.1666
$ws.Range("K107").Value = 1405.4667
$ws.Range("L107").Value = 1206.1666
$ws.Range("M107").Value = 514.5333000000001
$ws.Range("N107").Value = -5046.1666
$ws.Range("H134").Value = 9908.429
$ws.Range("I134").Value = 6086.3687
$ws.Range("J134").Value = 17977.223
$ws.Range("K134").Value = 18259.1061
$ws.Range("L134").Value = 53931.66900000001
$ws.Range("M134").Value = -15724.1061
$ws.Range("N134").Value = -59001.66900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2132.0625
$ws.Range("I31").Value = 1051
$ws.Range("K31").Value = 1051
$ws.Range("M31").Value = -756
$ws.Range("H33").Value = 11235
$ws.Range("I33").Value = 8898.75
$ws.Range("J33").Value = 17465
$ws.Range("K33").Value = 8898.75
$ws.Range("L33").Value = 17465
$ws.Range("M33").Value = -8519.75
$ws.Range("N33").Value = -18223
$ws.Range("H34").Value = 2132.0625
$ws.Range("I34").Value = 1051
$ws.Range("K34").Value = 1051
$ws.Range("M34").Value = -849
$ws.Range("H50").Value = 38745.25
$ws.Range("J50").Value = 44994
$ws.Range("L50").Value = 44994
$ws.Range("N50").Value = -46244
$ws.Range("H51").Value = 44994
$ws.Range("J51").Value = 44994
$ws.Range("L51").Value = 44994
$ws.Range("N51").Value = -46466
$ws.Range("H58").Value = 4107.269
$ws.Range("I58").Value = 2706.182
$ws.Range("K58").Value = 2706.182
$ws.Range("M58").Value = -2503.182
$ws.Range("H59").Value = 51842.965
$ws.Range("J59").Value = 51790.48
$ws.Range("L59").Value = 51790.48
$ws.Range("N59").Value = -54080.48
$ws.Range("H60").Value = 35724.57
$ws.Range("J60").Value = 34994
$ws.Range("L60").Value = 34994
$ws.Range("N60").Value = -36016
$ws.Range("H61").Value = 44994
$ws.Range("J61").Value = 44994
$ws.Range("L61").Value = 44994
$ws.Range("N61").Value = -45690
$ws.Range("H74").Value = 40795.2
$ws.Range("J74").Value = 44994
$ws.Range("L74").Value = 44994
$ws.Range("N74").Value = -46742
$ws.Range("H77").Value = 40795.2
$ws.Range("J77").Value = 44994
$ws.Range("L77").Value = 134982
$ws.Range("N77").Value = -143718
$ws.Range("H94").Value = 915.85
$ws.Range("J94").Value = 1060.7333
$ws.Range("L94").Value = 1060.7333
$ws.Range("N94").Value = -1962.7333
$ws.Range("H99").Value = 11315.125
$ws.Range("I99").Value = 7454.7144
$ws.Range("J99").Value = 12396.04
$ws.Range("K99").Value = 7454.7144
$ws.Range("L99").Value = 12396.04
$ws.Range("M99").Value = -5956.7144
$ws.Range("N99").Value = -15392.04
$ws.Range("H122").Value = 2148.1667
$ws.Range("I122").Value = 1778.1666
$ws.Range("K122").Value = 5334.4998
$ws.Range("M122").Value = -2884.4998
$ws.Range("H126").Value = 11315.125
$ws.Range("I126").Value = 7454.7144
$ws.Range("J126").Value = 12396.04
$ws.Range("K126").Value = 22364.1432
$ws.Range("L126").Value = 37188.12
$ws.Range("M126").Value = -19894.1432
$ws.Range("N126").Value = -42128.12
$ws.Range("H132").Value = 18320.05
$ws.Range("I132").Value = 10350.65
$ws.Range("J132").Value = 43486.58
$ws.Range("K132").Value = 31051.95
$ws.Range("L132").Value = 130459.74
$ws.Range("M132").Value = -28521.95
$ws.Range("N132").Value = -135519.74
$ws.Range("H134").Value = 4755.732
$ws.Range("I134").Value = 3922.261
$ws.Range("K134").Value = 11766.783
$ws.Range("M134").Value = -9231.782999999999
$ws.Range("H136").Value = 4107.269
$ws.Range("I136").Value = 2706.182
$ws.Range("K136").Value = 8118.545999999999
$ws.Range("M136").Value = -5568.545999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1097.6471
$ws.Range("I5").Value = 727.8823
$ws.Range("J5").Value = 1837.1765
$ws.Range("K5").Value = 2183.6469
$ws.Range("L5").Value = 5511.529500000001
$ws.Range("M5").Value = -2071.6469
$ws.Range("N5").Value = -5735.529500000001
$ws.Range("H34").Value = 2465.28
$ws.Range("J34").Value = 5285.273
$ws.Range("L34").Value = 15855.819
$ws.Range("N34").Value = -16023.819
$ws.Range("H40").Value = 93.71429000000001
$ws.Range("I40").Value = 90.25
$ws.Range("J40").Value = 98.333336
$ws.Range("K40").Value = 361
$ws.Range("L40").Value = 393.333344
$ws.Range("M40").Value = -292
$ws.Range("N40").Value = -531.333344
$ws.Range("H132").Value = 5904622.5
$ws.Range("I132").Value = 1631
$ws.Range("K132").Value = 14679
$ws.Range("M132").Value = -12149
$ws.Range("H135").Value = 1097.6471
$ws.Range("I135").Value = 727.8823
$ws.Range("J135").Value = 1837.1765
$ws.Range("K135").Value = 6550.9407
$ws.Range("L135").Value = 16534.5885
$ws.Range("M135").Value = -4015.9407
$ws.Range("N135").Value = -21604.5885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 31335
$ws.Range("J62").Value = 31335
$ws.Range("L62").Value = 31335
$ws.Range("N62").Value = -32707
$ws.Range("H65").Value = 31335
$ws.Range("J65").Value = 31335
$ws.Range("L65").Value = 94005
$ws.Range("N65").Value = -100869
$ws.Range("H80").Value = 6061.9614
$ws.Range("I80").Value = 3616.5334
$ws.Range("J80").Value = 9396.637000000001
$ws.Range("K80").Value = 3616.5334
$ws.Range("L80").Value = 9396.637000000001
$ws.Range("M80").Value = -2618.5334
$ws.Range("N80").Value = -11392.637
$ws.Range("H83").Value = 6061.9614
$ws.Range("I83").Value = 3616.5334
$ws.Range("J83").Value = 9396.637000000001
$ws.Range("K83").Value = 18082.667
$ws.Range("L83").Value = 46983.185
$ws.Range("M83").Value = -13090.667
$ws.Range("N83").Value = -56967.185
$ws.Range("H98").Value = 32108.6
$ws.Range("J98").Value = 32108.6
$ws.Range("L98").Value = 32108.6
$ws.Range("N98").Value = -38098.6
$ws.Range("H99").Value = 17251.7
$ws.Range("J99").Value = 46654
$ws.Range("L99").Value = 46654
$ws.Range("N99").Value = -51146
$ws.Range("H102").Value = 2244.682
$ws.Range("I102").Value = 2257.5264
$ws.Range("K102").Value = 2257.5264
$ws.Range("M102").Value = -635.5264000000002
$ws.Range("H113").Value = 114501.7
$ws.Range("I113").Value = 22904
$ws.Range("K113").Value = 22904
$ws.Range("M113").Value = -20734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8827.947
$ws.Range("I7").Value = 8858.1875
$ws.Range("J7").Value = 8666.666999999999
$ws.Range("K7").Value = 8858.1875
$ws.Range("L7").Value = 8666.666999999999
$ws.Range("M7").Value = -8746.1875
$ws.Range("N7").Value = -8890.666999999999
$ws.Range("H46").Value = 1261.1522
$ws.Range("I46").Value = 1122
$ws.Range("J46").Value = 1304.8857
$ws.Range("K46").Value = 1122
$ws.Range("L46").Value = 1304.8857
$ws.Range("M46").Value = -934
$ws.Range("N46").Value = -1680.8857
$ws.Range("H55").Value = 81.666664
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 3134.7
$ws.Range("I61").Value = 3004.25
$ws.Range("K61").Value = 3004.25
$ws.Range("M61").Value = -2802.25
$ws.Range("H82").Value = 2258.6
$ws.Range("I82").Value = 1630.091
$ws.Range("J82").Value = 3026.7778
$ws.Range("K82").Value = 1630.091
$ws.Range("L82").Value = 3026.7778
$ws.Range("M82").Value = -1269.091
$ws.Range("N82").Value = -3748.7778
$ws.Range("H85").Value = 2258.6
$ws.Range("I85").Value = 1630.091
$ws.Range("J85").Value = 3026.7778
$ws.Range("K85").Value = 1630.091
$ws.Range("L85").Value = 3026.7778
$ws.Range("M85").Value = -382.0909999999999
$ws.Range("N85").Value = -5522.7778
$ws.Range("H100").Value = 2025.7059
$ws.Range("I100").Value = 2162.6
$ws.Range("K100").Value = 2162.6
$ws.Range("M100").Value = -1621.6
$ws.Range("H110").Value = 53987.8
$ws.Range("J110").Value = 53987.8
$ws.Range("L110").Value = 53987.8
$ws.Range("N110").Value = -62167.8
$ws.Range("H113").Value = 3134.7
$ws.Range("I113").Value = 3004.25
$ws.Range("K113").Value = 3004.25
$ws.Range("M113").Value = -834.25
$ws.Range("H126").Value = 8827.947
$ws.Range("I126").Value = 8858.1875
$ws.Range("J126").Value = 8666.666999999999
$ws.Range("K126").Value = 26574.5625
$ws.Range("L126").Value = 26000.001
$ws.Range("M126").Value = -24104.5625
$ws.Range("N126").Value = -30940.001
$ws.Range("H132").Value = 3049.06
$ws.Range("I132").Value = 2675.923
$ws.Range("J132").Value = 4372
$ws.Range("K132").Value = 8027.768999999999
$ws.Range("L132").Value = 13116
$ws.Range("M132").Value = -5497.768999999999
$ws.Range("N132").Value = -18176
$ws.Range("H140").Value = 62339.375
$ws.Range("J140").Value = 62339.375
$ws.Range("L140").Value = 62339.375
$ws.Range("N140").Value = -72699.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 1003
$ws.Range("I22").Value = 1003
$ws.Range("K22").Value = 1003
$ws.Range("M22").Value = -710
$ws.Range("H69").Value = 33134
$ws.Range("J69").Value = 33134
$ws.Range("L69").Value = 33134
$ws.Range("N69").Value = -34632
$ws.Range("H72").Value = 33134
$ws.Range("J72").Value = 33134
$ws.Range("L72").Value = 99402
$ws.Range("N72").Value = -106890
$ws.Range("H75").Value = 36683.6
$ws.Range("J75").Value = 25879.75
$ws.Range("L75").Value = 25879.75
$ws.Range("N75").Value = -27751.75
$ws.Range("H78").Value = 36683.6
$ws.Range("J78").Value = 25879.75
$ws.Range("L78").Value = 77639.25
$ws.Range("N78").Value = -86999.25
$ws.Range("H81").Value = 2186.5454
$ws.Range("I81").Value = 1881.625
$ws.Range("K81").Value = 3763.25
$ws.Range("M81").Value = -2702.25
$ws.Range("H84").Value = 2186.5454
$ws.Range("I84").Value = 1881.625
$ws.Range("K84").Value = 18816.25
$ws.Range("M84").Value = -13512.25
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = None
$ws.Range("N95").ClearContents()
$ws.Range("H126").Value = 6051.607
$ws.Range("I126").Value = 6324.846
$ws.Range("J126").Value = 2499.5
$ws.Range("K126").Value = 18974.538
$ws.Range("L126").Value = 7498.5
$ws.Range("M126").Value = -16504.538
$ws.Range("N126").Value = -12438.5
$ws.Range("H132").Value = 14457.155
$ws.Range("I132").Value = 8358.727999999999
$ws.Range("J132").Value = 31227.834
$ws.Range("K132").Value = 25076.184
$ws.Range("L132").Value = 93683.50199999999
$ws.Range("M132").Value = -22546.184
$ws.Range("N132").Value = -98743.50199999999
$ws.Range("H136").Value = 1028.7916
$ws.Range("I136").Value = 1028.7916
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3086.3748
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = None
$ws.Range("N136").ClearContents()
